$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the "Date" property value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) ---
# (header row + the 5 data rows below it), including the column widths,
# so the "Spécification métier" mapping column now comes before "RIM Mapping".
$elements = $wb.Worksheets.Item("Elements")

# Only the rows whose AK/AL content actually differ need touching (rows
# that are blank in both columns must stay untouched so they remain blank
# cells rather than turning into empty-string text cells).
$rowsToSwap = @(1, 3, 5, 6)
foreach ($r in $rowsToSwap) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the two columns' widths to match their swapped content.
$elements.Columns.Item(37).ColumnWidth = 70.0703125
$elements.Columns.Item(38).ColumnWidth = 24.98046875
